# Auto commit at 2025-11-09 10:04:41.97
# Appends two new daily-charging-volume rows (2025-11-08) to Sheet1,
# one for 四方坪站充电量(kw) and one for 高岭站充电量(kw).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 45969

$row138 = @(587.23,1326.1529999999998,471.2,631.96,492.0680000000001,741.75599999999997,470.78200000000004,194.02499999999998,161.10899999999998,102.57999999999998,200.70499999999998,296.52600000000007,1016.7249999999999,1704.3670000000002,705.08799999999985,769.94799999999975,317.54199999999997,289.54700000000003,187.78799999999998,108.67999999999999,59.54,84.72,75.09,4.72)
$row139 = @(323.71000000000004,538.50800000000004,0,114.75800000000001,132.464,183.27700000000002,161.20299999999997,130.14399999999998,111.58599999999998,249.32199999999997,210.53100000000001,115.703,413.06900000000007,549.46100000000013,322.13300000000004,138.65700000000001,128.86800000000002,122.10900000000001,125.74299999999999,66.885999999999996,54.328000000000003,0.76600000000000001,25.75,0)

# Row 138: 四方坪站充电量(kw)
$ws.Cells.Item(138, 1).Value = $newDate
$ws.Cells.Item(138, 2).Value = "四方坪站充电量(kw)"
for ($i = 0; $i -lt $row138.Length; $i++) {
    $ws.Cells.Item(138, 3 + $i).Value = $row138[$i]
}

# Row 139: 高岭站充电量(kw)
$ws.Cells.Item(139, 1).Value = $newDate
$ws.Cells.Item(139, 2).Value = "高岭站充电量(kw)"
for ($i = 0; $i -lt $row139.Length; $i++) {
    $ws.Cells.Item(139, 3 + $i).Value = $row139[$i]
}

# Update the view: scroll position and active selection, matching the
# author's on-screen state after appending the new rows.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 110
$win.ScrollColumn = 1
$ws.Range("H142").Select() | Out-Null
